$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("L2").Value = 4245
$ws.Range("L3").Value = 4504
$ws.Range("L4").Value = 1114
$ws.Range("L5").Value = 258
$ws.Range("L6").Value = 3874
$ws.Range("L7").Value = 13995

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("L2").Value = 114
$ws.Range("L4").Value = 55
$ws.Range("L7").Value = 463
$ws.Range("L8").Value = 928
$ws.Range("L11").Value = 223
$ws.Range("L12").Value = 33
$ws.Range("L15").Value = 106
$ws.Range("L20").Value = 351
$ws.Range("L31").Value = 134
$ws.Range("L34").Value = 85
$ws.Range("L36").Value = 177
$ws.Range("L37").Value = 516
$ws.Range("L41").Value = 63
$ws.Range("L42").Value = 439
$ws.Range("L48").Value = 182
$ws.Range("L52").Value = 286
$ws.Range("L53").Value = 164
$ws.Range("L54").Value = 291
$ws.Range("L55").Value = 136
$ws.Range("L57").Value = 52
$ws.Range("K63").Value = 143
$ws.Range("L63").Value = 47
$ws.Range("L65").Value = 266
$ws.Range("L67").Value = 480
$ws.Range("L68").Value = 43
$ws.Range("L69").Value = 36
$ws.Range("L71").Value = 38
$ws.Range("L76").Value = 209
$ws.Range("K77").Value = 125
$ws.Range("L79").Value = 370
$ws.Range("L80").Value = 44
$ws.Range("L83").Value = 310
$ws.Range("L85").Value = 727
$ws.Range("L86").Value = 109
$ws.Range("L89").Value = 200
$ws.Range("L92").Value = 41
$ws.Range("L93").Value = 73
$ws.Range("L96").Value = 150
$ws.Range("L99").Value = 233
$ws.Range("L101").Value = 13995

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("L5").Value = 2
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 150

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("L2").Value = 154
$ws.Range("L7").Value = 463

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("L2").Value = 84
$ws.Range("L7").Value = 223

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("L2").Value = 55
$ws.Range("L7").Value = 200

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("L2").Value = 222
$ws.Range("L3").Value = 293
$ws.Range("L7").Value = 727

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("L3").Value = 93
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 286

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("L2").Value = 51
$ws.Range("L4").Value = 17
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("L2").Value = 268
$ws.Range("L3").Value = 309
$ws.Range("L5").Value = 33
$ws.Range("L6").Value = 252
$ws.Range("L7").Value = 928

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("L3").Value = 121
$ws.Range("L7").Value = 310

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("L3").Value = 180
$ws.Range("L7").Value = 516

$ws = $wb.Worksheets.Item('New City')
$ws.Range("L2").Value = 93
$ws.Range("L7").Value = 266

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("L3").Value = 100
$ws.Range("L7").Value = 233

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Range("L6").Value = 40
$ws.Range("L7").Value = 134

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("L2").Value = 137
$ws.Range("L3").Value = 184
$ws.Range("L7").Value = 480

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("L3").Value = 70
$ws.Range("L4").Value = 24
$ws.Range("L6").Value = 141
$ws.Range("L7").Value = 291

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("L3").Value = 43
$ws.Range("L7").Value = 182

$ws = $wb.Worksheets.Item('River North')
$ws.Range("L4").Value = 27
$ws.Range("L7").Value = 209

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("L3").Value = 22
$ws.Range("L4").Value = 7
$ws.Range("L7").Value = 63

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("L3").Value = 141
$ws.Range("L6").Value = 124
$ws.Range("L7").Value = 439

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("L6").Value = 39
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("L2").Value = 125
$ws.Range("L3").Value = 134
$ws.Range("L7").Value = 370

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("L3").Value = 109
$ws.Range("L6").Value = 98
$ws.Range("L7").Value = 351

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("L2").Value = 68
$ws.Range("L4").Value = 13
$ws.Range("L7").Value = 177

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("L6").Value = 24
$ws.Range("L7").Value = 73

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 85

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("L2").Value = 38
$ws.Range("L7").Value = 114

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 41

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range("L4").Value = 59
$ws.Range("L7").Value = 109

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 43

$ws = $wb.Worksheets.Item('Mckinley Park')
$ws.Range("L2").Value = 20
$ws.Range("L7").Value = 52

$ws = $wb.Worksheets.Item('Oakland')
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 125

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("L6").Value = 23
$ws.Range("L7").Value = 44

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range("L3").Value = 17
$ws.Range("L7").Value = 55

$ws = $wb.Worksheets.Item('Beverly')
$ws.Range("L3").Value = 8
$ws.Range("L7").Value = 33
